$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.709.20'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.278.26'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '''122.74'
$ws.Range('E5').Value = '  +7.38%  '
$ws.Range('D6').Value = '''265.32'
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('D7').Value = '''0.638'
$ws.Range('E7').Value = '  +1.88%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = '''0.622'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').Value = '''48.19'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').Value = '''0.0942'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '''8.99'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '''15.46'
$ws.Range('E14').Value = '  -2.34%  '
$ws.Range('D15').Value = '''0.897'
$ws.Range('E15').Value = '  +4.66%  '
$ws.Range('D16').Value = '2.624.99'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '2.279.90'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '43.678.46'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = '''6.99'
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '''72.39'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '''2.44'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '''235.48'
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').Value = '''9.53'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('D25').Value = '''2.87'
$ws.Range('E25').Value = '  -3.61%  '
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('D27').Value = '''11.80'
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').Value = '''42.00'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').Value = '''3.36'
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''171.88'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').Value = '''21.67'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('D33').Value = '''0.0909'
$ws.Range('E33').Value = '  -2.13%  '
$ws.Range('D34').Value = '''5.72'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '''0.129'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').Value = '''0.0379'
$ws.Range('E36').Value = '  +4.35%  '
$ws.Range('D37').Value = '''4.68'
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').Value = '''4.09'
$ws.Range('E38').Value = '  +7.09%  '
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').Value = '''2.52'
$ws.Range('E40').Value = '  +5.31%  '
$ws.Range('D41').Value = '''75.19'
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('D42').Value = '''13.84'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('D43').Value = '''0.238'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value = '''1.37'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '''5.76'
$ws.Range('E46').Value = '  -9.53%  '
$ws.Range('D47').Value = '''74.07'
$ws.Range('E47').Value = '  +38.30%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '''8.56'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').Value = '''1.26'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = '''0.100'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').Value = '''101.58'
$ws.Range('E51').Value = '  -0.85%  '
